$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A15").Value = "Website will crash if try to delete nil app"
$ws.Range("B15").Value = "resolved"

$ws.Range("A16").Select()
